# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-data refresh to the Leve profit tables
# across all 8 crafting-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ALC = $wb.Worksheets.Item("ALC")
# Row 15
$ALC.Range("H15").Value2 = 3876.7112
$ALC.Range("I15").Value2 = 3876.7112
$ALC.Range("K15").Value2 = 11630.1336
$ALC.Range("M15").Value2 = -11461.1336

# Row 20
$ALC.Range("H20").Value2 = 291.5
$ALC.Range("I20").Value2 = 291.5
$ALC.Range("K20").Value2 = 291.5
$ALC.Range("M20").Value2 = -61.5

# Row 35
$ALC.Range("H35").Value2 = 291.5
$ALC.Range("I35").Value2 = 291.5
$ALC.Range("K35").Value2 = 291.5
$ALC.Range("M35").Value2 = 87.5

# Row 113
$ALC.Range("H113").Value2 = 2550
$ALC.Range("I113").Value2 = 0
$ALC.Range("J113").Value2 = 2550
$ALC.Range("K113").Value2 = 0
$ALC.Range("L113").Value2 = 2550
$ALC.Range("M113").ClearContents()
$ALC.Range("N113").Value2 = -9058

# Row 125
$ALC.Range("H125").Value2 = 1199.2
$ALC.Range("J125").Value2 = 2000
$ALC.Range("L125").Value2 = 18000
$ALC.Range("N125").Value2 = -22920

# Row 132
$ALC.Range("H132").Value2 = 12078.947
$ALC.Range("I132").Value2 = 12699.667
$ALC.Range("J132").Value2 = 906
$ALC.Range("K132").Value2 = 38099.001
$ALC.Range("L132").Value2 = 2718
$ALC.Range("M132").Value2 = -35569.001
$ALC.Range("N132").Value2 = -7778

$ARM = $wb.Worksheets.Item("ARM")
# Row 2
$ARM.Range("H2").Value2 = 1250
$ARM.Range("I2").Value2 = 1250
$ARM.Range("J2").Value2 = 0
$ARM.Range("K2").Value2 = 1250
$ARM.Range("L2").Value2 = 0
$ARM.Range("M2").Value2 = -1137
$ARM.Range("N2").ClearContents()

# Row 32
$ARM.Range("H32").Value2 = 7693077.5
$ARM.Range("I32").Value2 = 834.1667
$ARM.Range("K32").Value2 = 834.1667
$ARM.Range("M32").Value2 = -547.1667

# Row 116
$ARM.Range("H116").Value2 = 1250
$ARM.Range("I116").Value2 = 1250
$ARM.Range("J116").Value2 = 0
$ARM.Range("K116").Value2 = 1250
$ARM.Range("L116").Value2 = 0
$ARM.Range("M116").Value2 = 1044
$ARM.Range("N116").ClearContents()

# Row 122
$ARM.Range("H122").Value2 = 0
$ARM.Range("I122").Value2 = 0
$ARM.Range("K122").Value2 = 0
$ARM.Range("M122").ClearContents()

# Row 132
$ARM.Range("H132").Value2 = 1628.8182
$ARM.Range("I132").Value2 = 1628.8182
$ARM.Range("K132").Value2 = 4886.4546
$ARM.Range("M132").Value2 = -2356.4546

$BSM = $wb.Worksheets.Item("BSM")
# Row 3
$BSM.Range("H3").Value2 = 1250
$BSM.Range("I3").Value2 = 1250
$BSM.Range("J3").Value2 = 0
$BSM.Range("K3").Value2 = 1250
$BSM.Range("L3").Value2 = 0
$BSM.Range("M3").Value2 = -1136
$BSM.Range("N3").ClearContents()

# Row 86
$BSM.Range("H86").Value2 = 5257.2856
$BSM.Range("I86").Value2 = 2900.2222
$BSM.Range("J86").Value2 = 9500
$BSM.Range("K86").Value2 = 2900.2222
$BSM.Range("L86").Value2 = 9500
$BSM.Range("M86").Value2 = -1777.2222
$BSM.Range("N86").Value2 = -11746

# Row 89
$BSM.Range("H89").Value2 = 5257.2856
$BSM.Range("I89").Value2 = 2900.2222
$BSM.Range("J89").Value2 = 9500
$BSM.Range("K89").Value2 = 14501.111
$BSM.Range("L89").Value2 = 47500
$BSM.Range("M89").Value2 = -8885.111000000001
$BSM.Range("N89").Value2 = -58732

# Row 134
$BSM.Range("H134").Value2 = 5150.706
$BSM.Range("I134").Value2 = 1974.0769
$BSM.Range("K134").Value2 = 5922.2307
$BSM.Range("M134").Value2 = -3387.2307

$CRP = $wb.Worksheets.Item("CRP")
# Row 6
$CRP.Range("H6").Value2 = 323.25
$CRP.Range("I6").Value2 = 323.25
$CRP.Range("K6").Value2 = 323.25
$CRP.Range("M6").Value2 = -210.25

# Row 33
$CRP.Range("H33").Value2 = 1191.0834
$CRP.Range("I33").Value2 = 666
$CRP.Range("K33").Value2 = 666
$CRP.Range("M33").Value2 = -287

# Row 99
$CRP.Range("H99").Value2 = 3594.6
$CRP.Range("I99").Value2 = 3594.6
$CRP.Range("K99").Value2 = 3594.6
$CRP.Range("M99").Value2 = -2096.6

# Row 126
$CRP.Range("H126").Value2 = 3594.6
$CRP.Range("I126").Value2 = 3594.6
$CRP.Range("K126").Value2 = 10783.8
$CRP.Range("M126").Value2 = -8313.799999999999

$CUL = $wb.Worksheets.Item("CUL")
# Row 86
$CUL.Range("H86").Value2 = 613.4286
$CUL.Range("I86").Value2 = 613.4286
$CUL.Range("J86").Value2 = 0
$CUL.Range("K86").Value2 = 1840.2858
$CUL.Range("L86").Value2 = 0
$CUL.Range("M86").Value2 = -654.2857999999999
$CUL.Range("N86").ClearContents()

# Row 89
$CUL.Range("H89").Value2 = 613.4286
$CUL.Range("I89").Value2 = 613.4286
$CUL.Range("J89").Value2 = 0
$CUL.Range("K89").Value2 = 5520.8574
$CUL.Range("L89").Value2 = 0
$CUL.Range("M89").Value2 = 407.1426000000001
$CUL.Range("N89").ClearContents()

# Row 122
$CUL.Range("H122").Value2 = 781.9231
$CUL.Range("I122").Value2 = 514.7143
$CUL.Range("J122").Value2 = 1093.6666
$CUL.Range("K122").Value2 = 4632.428699999999
$CUL.Range("L122").Value2 = 9842.999400000001
$CUL.Range("M122").Value2 = -2182.428699999999
$CUL.Range("N122").Value2 = -14742.9994

$GSM = $wb.Worksheets.Item("GSM")
# Row 101
$GSM.Range("H101").Value2 = 45997.5
$GSM.Range("J101").Value2 = 45997.5
$GSM.Range("L101").Value2 = 45997.5
$GSM.Range("N101").Value2 = -52487.5

# Row 122
$GSM.Range("H122").Value2 = 1120.0834
$GSM.Range("I122").Value2 = 1154.1
$GSM.Range("J122").Value2 = 950
$GSM.Range("K122").Value2 = 3462.3
$GSM.Range("L122").Value2 = 2850
$GSM.Range("M122").Value2 = -1012.3
$GSM.Range("N122").Value2 = -7750

$LTW = $wb.Worksheets.Item("LTW")
# Row 7
$LTW.Range("H7").Value2 = 3725.75
$LTW.Range("I7").Value2 = 2967.8333
$LTW.Range("K7").Value2 = 2967.8333
$LTW.Range("M7").Value2 = -2855.8333

# Row 40
$LTW.Range("H40").Value2 = 7810.4814
$LTW.Range("I40").Value2 = 7775.4
$LTW.Range("K40").Value2 = 7775.4
$LTW.Range("M40").Value2 = -7639.4

# Row 55
$LTW.Range("H55").Value2 = 1170.1052
$LTW.Range("I55").Value2 = 782.13336
$LTW.Range("J55").Value2 = 2625
$LTW.Range("K55").Value2 = 782.13336
$LTW.Range("L55").Value2 = 2625
$LTW.Range("M55").Value2 = -609.13336
$LTW.Range("N55").Value2 = -2971

# Row 122
$LTW.Range("H122").Value2 = 3016.5
$LTW.Range("I122").Value2 = 3018
$LTW.Range("J122").Value2 = 3000
$LTW.Range("K122").Value2 = 9054
$LTW.Range("L122").Value2 = 9000
$LTW.Range("M122").Value2 = -6604
$LTW.Range("N122").Value2 = -13900

# Row 126
$LTW.Range("H126").Value2 = 3725.75
$LTW.Range("I126").Value2 = 2967.8333
$LTW.Range("K126").Value2 = 8903.499899999999
$LTW.Range("M126").Value2 = -6433.499899999999

$WVR = $wb.Worksheets.Item("WVR")
# Row 104
$WVR.Range("H104").Value2 = 26967.25
$WVR.Range("J104").Value2 = 26967.25
$WVR.Range("L104").Value2 = 26967.25
$WVR.Range("N104").Value2 = -33955.25

# Row 110
$WVR.Range("H110").Value2 = 5000
$WVR.Range("J110").Value2 = 5000
$WVR.Range("L110").Value2 = 5000
$WVR.Range("N110").Value2 = -13180

# Row 116
$WVR.Range("H116").Value2 = 20000
$WVR.Range("J116").Value2 = 20000
$WVR.Range("L116").Value2 = 20000
$WVR.Range("N116").Value2 = -29178

# Row 126
$WVR.Range("H126").Value2 = 4328.278
$WVR.Range("I126").Value2 = 2180.111
$WVR.Range("K126").Value2 = 6540.333
$WVR.Range("M126").Value2 = -4070.333
